$wb = $excel.ActiveWorkbook

# --- summary_F3 sheet: update its remembered selection ---
$ws3 = $wb.Worksheets.Item("summary_F3")
$ws3.Activate()
$ws3.Range("B22").Select()

# --- rawdata_Clio sheet: add a new results row (row 4) ---
$ws1 = $wb.Worksheets.Item("rawdata_Clio")

$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(4, 3).Value = 2
$ws1.Cells.Item(4, 4).Value = 5
$ws1.Cells.Item(4, 5).Value = 1
$ws1.Cells.Item(4, 6).Value = 3
$ws1.Cells.Item(4, 7).Value = 4
$ws1.Cells.Item(4, 8).Value = 6

# Update the selection to reflect the newly entered row of data
$ws1.Range("A4:H4").Select()

# --- Make rawdata_Clio the active (selected) sheet/tab ---
$ws1.Activate()
